$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 34 (pushes the "Electical Team Items" block and
# everything below it down by one row), mirroring Excel's own behaviour
# when a user right-clicks a row header and chooses "Insert".
$ws.Rows(34).Insert()

# Populate the freshly inserted row with the new "Camera" line item.
$ws.Range("B34").Value = "Camera"
$ws.Range("C34").Value = "UDOO Autofocus Camera 5.0"
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 39
$ws.Range("F34").Formula = "=D34*E34"
$ws.Range("G34").Value = "http://shop.udoo.org/usa/accessories/autofocus-camera-5-0.html?___from_store=other&popup=no"

# Match the author's final selection recorded in the diff.
[void]$ws.Range("G34").Select()
